$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = " XXXXXXX"
$ws.Range("B2").Value = " ABC"
$ws.Range("C2").Value = 324234234

$ws.Range("A3").Value = " XXXXXXX"
$ws.Range("B3").Value = " ABC"
$ws.Range("C3").Value = 324234234

$ws.Range("A4").Value = " XXXXXXX"
$ws.Range("B4").Value = " ABC"
$ws.Range("C4").Value = 324234234

$ws.Range("A5").Value = " XXXXXXX"
$ws.Range("B5").Value = " ABC"
$ws.Range("C5").Value = 324234234
